$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.266.15'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '2.270.00'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '495.91'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.30'
$ws.Range('E6').Value = '  +2.39%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0967'
$ws.Range('E9').Value = '  +4.82%  '
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.330'
$ws.Range('E11').Value = '  +4.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.72'
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '2.680.68'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.27'
$ws.Range('E14').Value = '  +4.77%  '
$ws.Range('D15').Value = '54.209.49'
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('D17').Value = '2.272.91'
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.11'
$ws.Range('E18').Value = '  +4.68%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.11'
$ws.Range('E19').Value = '  +3.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '303.09'
$ws.Range('E20').Value = '  +2.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.46'
$ws.Range('E21').Value = '  +4.68%  '
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '61.77'
$ws.Range('E23').Value = '  -2.64%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '2.376.50'
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.150'
$ws.Range('E26').Value = '  +1.93%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.20'
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '170.16'
$ws.Range('E28').Value = '  +4.76%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.61'
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0682'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.87'
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.08'
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.70'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.895'
$ws.Range('E36').Value = '  +6.06%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.19'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.71'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.84'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.39'
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '126.04'
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.78'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0900'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.546'
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '238.17'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.371'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0205'
$ws.Range('E50').Value = '  +2.35%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.76'
$ws.Range('E51').Value = '  +0.93%  '
